# Adds three new worksheets with reference/lookup data tables:
#   e_sl_from_cbr, k_4, cht11_f_e1
$wb = $excel.ActiveWorkbook

# ---- e_sl_from_cbr ----
$ws1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws1.Name = "e_sl_from_cbr"
$ws1.Cells.Item(1, 1).Value = "cbr"
$ws1.Cells.Item(1, 2).Value = "e_sl"
$data1 = @(
    @(1.93514399541733, 9.1549295774647899),
    @(4.02619777864117, 17.746478873239401),
    @(4.8843996203230002, 20),
    @(5.8780160722749102, 22.112676056338),
    @(7.7287967790701497, 25.352112676056301),
    @(9.7613635875767795, 28.028169014084501),
    @(11.6528878421309, 30),
    @(16.6065615805622, 34.084507042253499),
    @(19.8245253140675, 36.197183098591502),
    @(25.240496996036899, 40),
    @(28.9426612471675, 43.521126760563298),
    @(28.9426612471675, 43.521126760563298),
    @(31.3692038539278, 47.183098591549196),
    @(33.1878425452485, 50),
    @(39.301180392810501, 59.8591549295774),
    @(46.540620354044997, 70),
    @(55.559105469558297, 80),
    @(65.793322465756702, 90),
    @(73.643236559027201, 96.338028169014095)
)
$r = 2
foreach ($row in $data1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r++
}
$ws1.Range("A2:B20").NumberFormat = "0.000"
$ws1.Columns.Item(1).ColumnWidth = 5.7109375
$ws1.Columns.Item(2).ColumnWidth = 5.7109375

$ws1.Sort.SortFields.Clear()
$ws1.Sort.SortFields.Add($ws1.Range("A2:A20")) | Out-Null
$ws1.Sort.SetRange($ws1.Range("A1:B20"))
$ws1.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$ws1.Sort.Apply()
$ws1.Range("D6").Select() | Out-Null

# ---- k_4 ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "k_4"
$ws2.Cells.Item(1, 1).Value = "f_c"
$ws2.Cells.Item(1, 2).Value = "k_4"
$data2 = @(
    @(20, 1.03),
    @(25, 1.07),
    @(32, 1.1100000000000001),
    @(40, 1.1599999999999999),
    @(50, 1.2)
)
$r = 2
foreach ($row in $data2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r++
}
$ws2.Range("B7").Select() | Out-Null

# ---- cht11_f_e1 ----
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "cht11_f_e1"
$ws3.Cells.Item(1, 1).Value = "e_ss"
$ws3.Cells.Item(1, 2).Value = "f_e1"
$data3 = @(
    @(3.55202295615514, 0.80252164798861203),
    @(4.0979411868127604, 0.84037536437590898),
    @(4.6639367133060396, 0.86044069237536502),
    @(5.2137894807093197, 0.88093782451179403),
    @(5.9008469340084702, 0.900423015524627),
    @(6.6294605107390501, 0.92000088103853195),
    @(7.4656095746560904, 0.94059402495084998),
    @(8.3762828337628292, 0.96007960990094598),
    @(9.3221406289356903, 0.98026345054557695),
    @(10.491698223964599, 1.00121285220791),
    @(11.602675678526699, 1.01947640254954),
    @(12.8302345976005, 1.03994390827682),
    @(14.1738912965008, 1.05954878750117),
    @(15.4591459683169, 1.0792573782549699),
    @(17.0361328370279, 1.1003706810990801),
    @(18.5842059897636, 1.1201938652938199),
    @(20.200594921304099, 1.1396423556205999),
    @(22.120966173590599, 1.15994853787513),
    @(23.772524567090301, 1.18006091450031),
    @(25.636720725665398, 1.1999729500595999),
    @(27.669363479074502, 1.2202018081342201),
    @(29.648768529821002, 1.2397373849434801),
    @(31.635083118982401, 1.2593929115506699),
    @(33.762007859842299, 1.2793849258935199),
    @(36.113371611133701, 1.29941416502408),
    @(38.398701253034602, 1.3201299541288001),
    @(41.104137049930202, 1.33962273510456),
    @(43.954759141928498, 1.3599986954910199),
    @(47.109573971095699, 1.3808399529611299),
    @(49.9845754116104, 1.39850544174004),
    @(52.404594340712599, 1.40816799070605),
    @(54.972472133057998, 1.4170305701721),
    @(57.921419995880797, 1.4253392595021901),
    @(59.983247183165801, 1.43015740135247),
    @(62.7062706270627, 1.43536491373688),
    @(64.931984065986498, 1.43867893742102),
    @(67.491749174917501, 1.44205845734274),
    @(69.880720948807195, 1.4439510911092299),
    @(72.607260726072596, 1.4460544258018599),
    @(74.911966082453006, 1.44839006448015),
    @(78.045521447121899, 1.4529323234185401),
    @(79.852260974448697, 1.4566208236529701),
    @(82.829241828292496, 1.46382206392856),
    @(84.973577266402501, 1.46928475981672),
    @(87.9911952926103, 1.4797825983187001),
    @(89.922314149223098, 1.48645467856604),
    @(92.396682590633404, 1.4959004459013601),
    @(94.623526279568594, 1.50681436725003),
    @(96.457401774097804, 1.51713196387897),
    @(98.243635805965695, 1.52648449421196),
    @(99.935417405151398, 1.5369989805425399),
    @(102.18866635521999, 1.5512042206163601),
    @(104.755665064223, 1.56506299368914),
    @(107.43054613680501, 1.5798744341909401),
    @(109.784448764511, 1.5905352581952401),
    @(112.010413370104, 1.5997667961292601),
    @(114.94134002441299, 1.61121620340967),
    @(117.326732673267, 1.6201288018023301),
    @(119.76355169763499, 1.6286496543797599),
    @(122.277227722772, 1.6355783482539801),
    @(124.712288580456, 1.6427495511615),
    @(127.39273927392701, 1.64900545443765),
    @(129.99017938323499, 1.6544751231952),
    @(132.343234323432, 1.6590657793324199),
    @(135.02142451687999, 1.66333164016303),
    @(137.62376237623701, 1.6664285290804499),
    @(139.970161399701, 1.6692862539473501),
    @(142.40924092409199, 1.6711011146024699),
    @(144.945303543655, 1.6722251931623),
    @(147.524752475247, 1.67374977767237),
    @(150.092178167588, 1.6746372247141801)
)
$r = 2
foreach ($row in $data3) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $r++
}
$ws3.Sort.SortFields.Clear()
$ws3.Sort.SortFields.Add($ws3.Range("A2:A72")) | Out-Null
$ws3.Sort.SetRange($ws3.Range("A1:B72"))
$ws3.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$ws3.Sort.Apply()
$ws3.Range("G12").Select() | Out-Null

